$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stations and Times")

# Set the new values for column D (rows 23-25)
$ws.Range("D23").Value = 14.9964190582785
$ws.Range("D24").Value = 14.481906897791999
$ws.Range("D25").Value = 14.2431107560074

# Set the new values for column H (rows 23-25)
$ws.Range("H23").Value = 16.921301270796899
$ws.Range("H24").Value = 14.734212145307801
$ws.Range("H25").Value = 14.5142202430478

# Apply number format (2 decimal places) to these cells
$ws.Range("D23:D25").NumberFormat = "0.00"
$ws.Range("H23:H25").NumberFormat = "0.00"

# D column is right-aligned already by style but ensure explicit right alignment matches diff (xf 13 has horizontal=right)
$ws.Range("D23:D25").HorizontalAlignment = -4152  # xlRight

# Update selection to match diff: H23:H25 with active cell H23
$ws.Range("H23:H25").Select()
